# Apply updates to the "Metadata" worksheet (first sheet) to bump the
# ValueSet version/date/publisher/jurisdiction info, matching the new
# Alvearie FHIR IG publication metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Old row 10 "Contact / No display for ContactDetail" becomes
# "Jurisdiction / United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Old row 11 was a duplicate "Contact / No display for ContactDetail" row
# that is removed entirely; everything below shifts up by one row.
$ws.Range("A11:B11").EntireRow.Delete()
